$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing A13 value (tiny precision correction)
$ws.Range("A13").Value = 45877.50019596065

# Add new row 14 data
$ws.Range("A14").Value = 45877.54186516787
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 32
$ws.Range("D14").Value = 18.44
$ws.Range("E14").Value = 79.11
$ws.Range("F14").Value = 622.27
$ws.Range("G14").Value = 12.2
$ws.Range("H14").Value = "ESE"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "13:00:17"

# Copy style from A13 to A14 (date format with border, etc.)
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
